$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 96 (shifts the old rows 96-196 down to 98-198)
$ws.Range("A96:A97").EntireRow.Insert()

# New row 96: Comercializadora del Agro de Limarí, Poroto verde, Magnum
$ws.Cells.Item(96, 1).Value = 2
$ws.Cells.Item(96, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(96, 3).Value = "Coquimbo"
$ws.Cells.Item(96, 4).Value = 44679
$ws.Cells.Item(96, 5).Value = 4
$ws.Cells.Item(96, 6).Value = 100112031
$ws.Cells.Item(96, 7).Value = "Poroto verde"
$ws.Cells.Item(96, 8).Value = "Magnum"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 500
$ws.Cells.Item(96, 11).Value = 13000
$ws.Cells.Item(96, 12).Value = 14000
$ws.Cells.Item(96, 13).Value = 13500
$ws.Cells.Item(96, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(96, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(96, 16).Value = 540
$ws.Cells.Item(96, 17).Value = 25
$ws.Cells.Item(96, 18).Value = "Hortaliza"

# New row 97: Comercializadora del Agro de Limarí, Poroto verde, Sin especificar
$ws.Cells.Item(97, 1).Value = 2
$ws.Cells.Item(97, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(97, 3).Value = "Coquimbo"
$ws.Cells.Item(97, 4).Value = 44679
$ws.Cells.Item(97, 5).Value = 4
$ws.Cells.Item(97, 6).Value = 100112031
$ws.Cells.Item(97, 7).Value = "Poroto verde"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 360
$ws.Cells.Item(97, 11).Value = 19000
$ws.Cells.Item(97, 12).Value = 20000
$ws.Cells.Item(97, 13).Value = 19500
$ws.Cells.Item(97, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(97, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(97, 16).Value = 780
$ws.Cells.Item(97, 17).Value = 25
$ws.Cells.Item(97, 18).Value = "Hortaliza"
